$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update requirement text / values for rows 2-11
$ws.Range("B2").Value = "500 users can use the applications at the same time"

$ws.Range("B5").Value = "The system will scale to 100 GB of data"
$ws.Range("C5").Value = 5

$ws.Range("B6").Value = "The systems components will be independent"
$ws.Range("C6").Value = 3

$ws.Range("B7").Value = "There will be not loss of data due to connection issues"
$ws.Range("C7").Value = 5

$ws.Range("B8").Value = "No critical bugs in the system"
$ws.Range("D8").Value = 5

$ws.Range("B9").Value = "User interfaces will be represented in Hebrew"
$ws.Range("C9").Value = 5
$ws.Range("D9").Value = 1

$ws.Range("B10").Value = "Detailed logs will be written for each component"
$ws.Range("C10").Value = 5

$ws.Range("B11").Value = "Data will never be permanently deleted from the Database"
$ws.Range("C11").Value = 4
$ws.Range("D11").Value = 1

# Remove the now-obsolete extra requirement rows (12-15) entirely
$ws.Range("A12:G15").Clear()

# Clear the leftover row-number markers that used to sit in rows 16-17
# (keep their existing formatting, only drop the values)
$ws.Range("A16:G17").ClearContents()

# Fix the header typo: "Requirment" -> "Requirement"
$ws.Range("B1").Value = "Requirement"

# Update the active cell selection
$ws.Range("C19").Select()
